{"js": "const replacements = [\n  [\"68\u00d762=\", \"23\u00d732=\"],\n  [\"12\u00d794=\", \"41\u00d728=\"],\n  [\"43\u00d723=\", \"23\u00d783=\"],\n  [\"19\u00d740=\", \"41\u00d776=\"],\n  [\"60\u00d716=\", \"25\u00d721=\"],\n  [\"74\u00d737=\", \"68\u00d714=\"],\n  [\"34\u00d779=\", \"19\u00d797=\"],\n  [\"48\u00d785=\", \"77\u00d795=\"],\n  [\"75\u00d741=\", \"49\u00d715=\"],\n  [\"31\u00d797=\", \"75\u00d774=\"],\n  [\"43\u00d775=\", \"59\u00d765=\"],\n  [\"90\u00d773=\", \"42\u00d775=\"],\n  [\"77\u00d779=\", \"53\u00d768=\"],\n  [\"68\u00d755=\", \"72\u00d794=\"],\n  [\"81\u00d738=\", \"98\u00d748=\"],\n  [\"89\u00d791=\", \"99\u00d748=\"],\n  [\"63\u00d715=\", \"50\u00d756=\"],\n  [\"87\u00d787=\", \"49\u00d767=\"],\n  [\"89\u00d761=\", \"23\u00d724=\"],\n  [\"94\u00d770=\", \"64\u00d725=\"],\n  [\"79\u00d723=\", \"77\u00d732=\"],\n  [\"64\u00d758=\", \"71\u00d778=\"],\n  [\"74\u00d745=\", \"65\u00d796=\"],\n  [\"67\u00d788=\", \"59\u00d727=\"],\n  [\"14\u00d766=\", \"16\u00d791=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"68\u00d762=\", \"23\u00d732=\"),\n    @(\"12\u00d794=\", \"41\u00d728=\"),\n    @(\"43\u00d723=\", \"23\u00d783=\"),\n    @(\"19\u00d740=\", \"41\u00d776=\"),\n    @(\"60\u00d716=\", \"25\u00d721=\"),\n    @(\"74\u00d737=\", \"68\u00d714=\"),\n    @(\"34\u00d779=\", \"19\u00d797=\"),\n    @(\"48\u00d785=\", \"77\u00d795=\"),\n    @(\"75\u00d741=\", \"49\u00d715=\"),\n    @(\"31\u00d797=\", \"75\u00d774=\"),\n    @(\"43\u00d775=\", \"59\u00d765=\"),\n    @(\"90\u00d773=\", \"42\u00d775=\"),\n    @(\"77\u00d779=\", \"53\u00d768=\"),\n    @(\"68\u00d755=\", \"72\u00d794=\"),\n    @(\"81\u00d738=\", \"98\u00d748=\"),\n    @(\"89\u00d791=\", \"99\u00d748=\"),\n    @(\"63\u00d715=\", \"50\u00d756=\"),\n    @(\"87\u00d787=\", \"49\u00d767=\"),\n    @(\"89\u00d761=\", \"23\u00d724=\"),\n    @(\"94\u00d770=\", \"64\u00d725=\"),\n    @(\"79\u00d723=\", \"77\u00d732=\"),\n    @(\"64\u00d758=\", \"71\u00d778=\"),\n    @(\"74\u00d745=\", \"65\u00d796=\"),\n    @(\"67\u00d788=\", \"59\u00d727=\"),\n    @(\"14\u00d766=\", \"16\u00d791=\"),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
